# Update the two header cells (A1/B1) with new "relation" / "count" labels,
# replacing the old single merged-looking description string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "relation"
$ws.Range("B1").Value = "count"

# Widen column A so the (now much shorter) header text still reads well
# against the long relation strings beneath it.
$ws.Columns("A").ColumnWidth = 73.25

# Leave the active selection on the new "count" header cell.
$ws.Range("B1").Select() | Out-Null
